# Applies the "10강. css 속성 1 ( url )" edit described by the diff:
#  1. Paragraph "497F8274" (empty, style a4) loses its pPr/rPr rFonts hint.
#  2. Paragraph "38512251" ("2202-02-04") merges the split "202-02-0"+"4" runs.
#  3. Paragraph "1B9E3E64" ("10강. css 속성 1") merges the split " "+"속성"+" " runs.
#  4. Paragraph "7F912B51" ("em : ...") loses its pPr/rPr rFonts hint.
#  5. Four new paragraphs are appended at the end of the body (an empty
#     paragraph, a repeat of the date line, a repeat of the "10강." line,
#     and a brand-new "url : Background-img ..." list paragraph).

$d = $word.ActiveDocument
$wNs = 'http://schemas.openxmlformats.org/wordprocessingml/2006/main'

# --- 1. Blank "a4" paragraph right before "2202-02-04": drop the stray rPr.
$p157 = $d.Paragraphs.Item(157)
$xml157 = '<w:p xmlns:w="' + $wNs + '"><w:pPr><w:pStyle w:val="a4"/></w:pPr></w:p>'
$null = $p157.Range.InsertXML($xml157)

# --- 2. "2202-02-04": merge the "202-02-0" / "4" runs into a single run.
$p158 = $d.Paragraphs.Item(158)
$null = $p158.Range.Find.Execute("202-02-04", $false, $false, $false, $false, $false, $true, 1, $false, "202-02-04", 2)

# --- 3. "10강. css 속성 1": merge " " / "속성" / " " into a single run.
$p159 = $d.Paragraphs.Item(159)
$null = $p159.Range.Find.Execute(" 속성 ", $false, $false, $false, $false, $false, $true, 1, $false, " 속성 ", 2)

# --- 4. "em : ..." list paragraph: drop the stray rPr (keep pStyle/numPr/ind).
$p163 = $d.Paragraphs.Item(163)
$xml163 = '<w:p xmlns:w="' + $wNs + '">' + `
    '<w:pPr><w:pStyle w:val="a4"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr><w:ind w:leftChars="0"/></w:pPr>' + `
    '<w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/>' + `
    '<w:r><w:t>em</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> :</w:t></w:r><w:proofErr w:type="gramEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> 1.0</w:t></w:r>' + `
    '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>을 기준으로 올라가면 배로 확대 내려가면 배로 축소이다.</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> (2.0</w:t></w:r>' + `
    '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>은 두배)</w:t></w:r>' + `
    '</w:p>'
$null = $p163.Range.InsertXML($xml163)

# --- 5. Append the four new trailing paragraphs after the last paragraph.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$endRange = $d.Range($lastPara.Range.End, $lastPara.Range.End)

$newBlock = `
    '<w:p xmlns:w="' + $wNs + '"/>' + `
    '<w:p xmlns:w="' + $wNs + '">' + `
        '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>2</w:t></w:r>' + `
        '<w:r><w:t>202-02-04</w:t></w:r>' + `
    '</w:p>' + `
    '<w:p xmlns:w="' + $wNs + '">' + `
        '<w:r><w:t>10</w:t></w:r>' + `
        '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>강.</w:t></w:r>' + `
        '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
        '<w:proofErr w:type="spellStart"/><w:r><w:t>css</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
        '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> 속성 </w:t></w:r>' + `
        '<w:r><w:t>1</w:t></w:r>' + `
    '</w:p>' + `
    '<w:p xmlns:w="' + $wNs + '">' + `
        '<w:pPr><w:pStyle w:val="a4"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:ind w:leftChars="0"/><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr>' + `
        '<w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/>' + `
        '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>u</w:t></w:r>' + `
        '<w:r><w:t>rl</w:t></w:r>' + `
        '<w:proofErr w:type="spellEnd"/>' + `
        '<w:r><w:t xml:space="preserve"> :</w:t></w:r>' + `
        '<w:proofErr w:type="gramEnd"/>' + `
        '<w:r><w:t xml:space="preserve"> Backgroun</w:t></w:r>' + `
        '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>d</w:t></w:r>' + `
        '<w:r><w:t>-</w:t></w:r>' + `
        '<w:proofErr w:type="spellStart"/><w:r><w:t>img</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
        '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
        '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>속성의 속성값으로 많이 사용된다.</w:t></w:r>' + `
        '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
        '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>이 경우 배경 이미지의 경로를 나타낸다.</w:t></w:r>' + `
    '</w:p>'

$null = $endRange.InsertXML($newBlock)

Write-Host "Done. Paragraph count:" $d.Paragraphs.Count
